$d = $word.ActiveDocument

# Phase 1: replace each original cell text with a unique placeholder token
# to avoid any cross-matching between old/new values that collide.
$d.Content.Find.Execute("90÷9=10, 0", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_00@@", 2) | Out-Null
$d.Content.Find.Execute("29÷4=7, 1", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_01@@", 2) | Out-Null
$d.Content.Find.Execute("42÷9=4, 6", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_02@@", 2) | Out-Null
$d.Content.Find.Execute("40÷6=6, 4", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_03@@", 2) | Out-Null
$d.Content.Find.Execute("52÷5=10, 2", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_04@@", 2) | Out-Null
$d.Content.Find.Execute("52÷8=6, 4", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_05@@", 2) | Out-Null
$d.Content.Find.Execute("50÷2=25, 0", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_06@@", 2) | Out-Null
$d.Content.Find.Execute("85÷3=28, 1", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_07@@", 2) | Out-Null
$d.Content.Find.Execute("21÷7=3, 0", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_08@@", 2) | Out-Null
$d.Content.Find.Execute("88÷5=17, 3", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_09@@", 2) | Out-Null
$d.Content.Find.Execute("94÷5=18, 4", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_10@@", 2) | Out-Null
$d.Content.Find.Execute("54÷5=10, 4", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_11@@", 2) | Out-Null
$d.Content.Find.Execute("62÷7=8, 6", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_12@@", 2) | Out-Null
$d.Content.Find.Execute("65÷8=8, 1", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_13@@", 2) | Out-Null
$d.Content.Find.Execute("85÷5=17, 0", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_14@@", 2) | Out-Null
$d.Content.Find.Execute("37÷4=9, 1", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_15@@", 2) | Out-Null
$d.Content.Find.Execute("11÷8=1, 3", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_16@@", 2) | Out-Null
$d.Content.Find.Execute("76÷6=12, 4", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_17@@", 2) | Out-Null
$d.Content.Find.Execute("30÷9=3, 3", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_18@@", 2) | Out-Null
$d.Content.Find.Execute("54÷4=13, 2", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_19@@", 2) | Out-Null
$d.Content.Find.Execute("73÷9=8, 1", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_20@@", 2) | Out-Null
$d.Content.Find.Execute("14÷9=1, 5", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_21@@", 2) | Out-Null
$d.Content.Find.Execute("99÷5=19, 4", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_22@@", 2) | Out-Null
$d.Content.Find.Execute("34÷4=8, 2", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_23@@", 2) | Out-Null
$d.Content.Find.Execute("20÷8=2, 4", $true, $false, $false, $false, $false, $true, 1, $false, "@@PLACEHOLDER_24@@", 2) | Out-Null

# Phase 2: replace placeholders with final new text
$d.Content.Find.Execute("@@PLACEHOLDER_00@@", $true, $false, $false, $false, $false, $true, 1, $false, "87÷8=10, 7", 2) | Out-Null
$d.Content.Find.Execute("@@PLACEHOLDER_01@@", $true, $false, $false, $false, $false, $true, 1, $false, "93÷9=10, 3", 2) | Out-Null
$d.Content.Find.Execute("@@PLACEHOLDER_02@@", $true, $false, $false, $false, $false, $true, 1, $false, "25÷8=3, 1", 2) | Out-Null
$d.Content.Find.Execute("@@PLACEHOLDER_03@@", $true, $false, $false, $false, $false, $true, 1, $false, "72÷7=10, 2", 2) | Out-Null
$d.Content.Find.Execute("@@PLACEHOLDER_04@@", $true, $false, $false, $false, $false, $true, 1, $false, "33÷4=8, 1", 2) | Out-Null
$d.Content.Find.Execute("@@PLACEHOLDER_05@@", $true, $false, $false, $false, $false, $true, 1, $false, "42÷5=8, 2", 2) | Out-Null
$d.Content.Find.Execute("@@PLACEHOLDER_06@@", $true, $false, $false, $false, $false, $true, 1, $false, "64÷3=21, 1", 2) | Out-Null
$d.Content.Find.Execute("@@PLACEHOLDER_07@@", $true, $false, $false, $false, $false, $true, 1, $false, "97÷7=13, 6", 2) | Out-Null
$d.Content.Find.Execute("@@PLACEHOLDER_08@@", $true, $false, $false, $false, $false, $true, 1, $false, "66÷7=9, 3", 2) | Out-Null
$d.Content.Find.Execute("@@PLACEHOLDER_09@@", $true, $false, $false, $false, $false, $true, 1, $false, "64÷3=21, 1", 2) | Out-Null
$d.Content.Find.Execute("@@PLACEHOLDER_10@@", $true, $false, $false, $false, $false, $true, 1, $false, "80÷9=8, 8", 2) | Out-Null
$d.Content.Find.Execute("@@PLACEHOLDER_11@@", $true, $false, $false, $false, $false, $true, 1, $false, "57÷4=14, 1", 2) | Out-Null
$d.Content.Find.Execute("@@PLACEHOLDER_12@@", $true, $false, $false, $false, $false, $true, 1, $false, "97÷7=13, 6", 2) | Out-Null
$d.Content.Find.Execute("@@PLACEHOLDER_13@@", $true, $false, $false, $false, $false, $true, 1, $false, "37÷4=9, 1", 2) | Out-Null
$d.Content.Find.Execute("@@PLACEHOLDER_14@@", $true, $false, $false, $false, $false, $true, 1, $false, "28÷4=7, 0", 2) | Out-Null
$d.Content.Find.Execute("@@PLACEHOLDER_15@@", $true, $false, $false, $false, $false, $true, 1, $false, "12÷9=1, 3", 2) | Out-Null
$d.Content.Find.Execute("@@PLACEHOLDER_16@@", $true, $false, $false, $false, $false, $true, 1, $false, "31÷8=3, 7", 2) | Out-Null
$d.Content.Find.Execute("@@PLACEHOLDER_17@@", $true, $false, $false, $false, $false, $true, 1, $false, "18÷6=3, 0", 2) | Out-Null
$d.Content.Find.Execute("@@PLACEHOLDER_18@@", $true, $false, $false, $false, $false, $true, 1, $false, "14÷2=7, 0", 2) | Out-Null
$d.Content.Find.Execute("@@PLACEHOLDER_19@@", $true, $false, $false, $false, $false, $true, 1, $false, "29÷8=3, 5", 2) | Out-Null
$d.Content.Find.Execute("@@PLACEHOLDER_20@@", $true, $false, $false, $false, $false, $true, 1, $false, "47÷5=9, 2", 2) | Out-Null
$d.Content.Find.Execute("@@PLACEHOLDER_21@@", $true, $false, $false, $false, $false, $true, 1, $false, "59÷6=9, 5", 2) | Out-Null
$d.Content.Find.Execute("@@PLACEHOLDER_22@@", $true, $false, $false, $false, $false, $true, 1, $false, "72÷9=8, 0", 2) | Out-Null
$d.Content.Find.Execute("@@PLACEHOLDER_23@@", $true, $false, $false, $false, $false, $true, 1, $false, "50÷2=25, 0", 2) | Out-Null
$d.Content.Find.Execute("@@PLACEHOLDER_24@@", $true, $false, $false, $false, $false, $true, 1, $false, "51÷3=17, 0", 2) | Out-Null

Write-Output "Done replacing 25 cell values."
